$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries: Kenia now listed before Guinea-Bisau (row 102/103 swap) ---
# Row 102 becomes Kenia's data, row 103 becomes Guinea-Bisau's data.
$ws.Range("A102").Value = "Kenia"
$ws.Range("A103").Value = "Guinea-Bisau"

# --- Update the "Datos actualizados" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 14:35"

# --- Update numeric stats per row ---
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 1621669
$ws.Range("C4").Value = 767
$ws.Range("E4").Value = 1143048
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 96377

# Row 14
$ws.Range("B14").Value = 119574
$ws.Range("C14").Value = 1348
$ws.Range("E14").Value = 66948

# Row 23
$ws.Range("B23").Value = 44888
$ws.Range("C23").Value = 188
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 5788

# Row 31
$ws.Range("B31").Value = 30200
$ws.Range("C31").Value = 288
$ws.Range("D31").Value = 7590
$ws.Range("E31").Value = 21321
$ws.Range("G31").Value = 12
$ws.Range("H31").Value = 1289

# Row 47
$ws.Range("B47").Value = 11230
$ws.Range("C47").Value = 48
$ws.Range("D47").Value = 9764
$ws.Range("E47").Value = 905

# Row 57
$ws.Range("D57").Value = 3885
$ws.Range("E57").Value = 3677

# Row 85
$ws.Range("B85").Value = 2243
$ws.Range("C85").Value = 6
$ws.Range("D85").Value = 2011
$ws.Range("E85").Value = 133
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 99

# Row 102 (now Kenia) - values previously on row 103 (Kenia's old data), plus updates
$ws.Range("B102").Value = 1161
$ws.Range("C102").Value = 52
$ws.Range("D102").Value = 375
$ws.Range("E102").Value = 736
$ws.Range("H102").Value = 50

# Row 103 (now Guinea-Bisau) - values previously on row 102 (Guinea-Bisau's old data)
$ws.Range("B103").Value = 1109
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 42
$ws.Range("E103").Value = 1061
$ws.Range("H103").Value = 6
